$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 817.5077
$ws.Range("J17").Value = 817.5077
$ws.Range("L17").Value = 2452.5231
$ws.Range("N17").Value = -2788.5231

$ws.Range("H33").Value = 119
$ws.Range("I33").Value = 123.75
$ws.Range("J33").Value = 100
$ws.Range("K33").Value = 123.75
$ws.Range("L33").Value = 100
$ws.Range("M33").Value = 105.25
$ws.Range("N33").Value = -558

$ws.Range("H129").Value = 892.3472
$ws.Range("J129").Value = 874.9857
$ws.Range("L129").Value = 2624.9571
$ws.Range("N129").Value = -12624.9571

$ws.Range("H132").Value = 1479.75
$ws.Range("I132").Value = 1287
$ws.Range("J132").Value = 2250.75
$ws.Range("K132").Value = 3861
$ws.Range("L132").Value = 6752.25
$ws.Range("M132").Value = -1331
$ws.Range("N132").Value = -11812.25

$ws.Range("H138").Value = 2816
$ws.Range("I138").Value = 2486.1482
$ws.Range("J138").Value = 3145.8518
$ws.Range("K138").Value = 7458.444600000001
$ws.Range("L138").Value = 9437.555399999999
$ws.Range("M138").Value = -2318.444600000001
$ws.Range("N138").Value = -19717.5554


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3097.2134
$ws.Range("I32").Value = 2318.9092
$ws.Range("K32").Value = 2318.9092
$ws.Range("M32").Value = -2031.9092

$ws.Range("H61").Value = 4810.7144
$ws.Range("I61").Value = 3369.4285
$ws.Range("K61").Value = 3369.4285
$ws.Range("M61").Value = -3157.4285

$ws.Range("H74").Value = 1771.0731
$ws.Range("I74").Value = 1705.25
$ws.Range("K74").Value = 1705.25
$ws.Range("M74").Value = -831.25

$ws.Range("H77").Value = 1771.0731
$ws.Range("I77").Value = 1705.25
$ws.Range("K77").Value = 8526.25
$ws.Range("M77").Value = -4158.25

$ws.Range("H102").Value = 1840.4615
$ws.Range("I102").Value = 1538.8182
$ws.Range("K102").Value = 1538.8182
$ws.Range("M102").Value = 83.18180000000007

$ws.Range("H132").Value = 2346.8
$ws.Range("I132").Value = 2044.8462
$ws.Range("J132").Value = 2673.9167
$ws.Range("K132").Value = 6134.5386
$ws.Range("L132").Value = 8021.750100000001
$ws.Range("M132").Value = -3604.5386
$ws.Range("N132").Value = -13081.7501

$ws.Range("H136").Value = 4810.7144
$ws.Range("I136").Value = 3369.4285
$ws.Range("K136").Value = 10108.2855
$ws.Range("M136").Value = -7558.2855


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2400.611
$ws.Range("I134").Value = 2549.0967
$ws.Range("J134").Value = 1480
$ws.Range("K134").Value = 7647.2901
$ws.Range("L134").Value = 4440
$ws.Range("M134").Value = -5112.2901
$ws.Range("N134").Value = -9510


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2072000
$ws.Range("I58").Value = 2899858.5
$ws.Range("J58").Value = 2353.8333
$ws.Range("K58").Value = 2899858.5
$ws.Range("L58").Value = 2353.8333
$ws.Range("M58").Value = -2899655.5
$ws.Range("N58").Value = -2759.8333

$ws.Range("H86").Value = 2052.5454
$ws.Range("J86").Value = 2432.5
$ws.Range("L86").Value = 2432.5
$ws.Range("N86").Value = -4678.5

$ws.Range("H89").Value = 2052.5454
$ws.Range("J89").Value = 2432.5
$ws.Range("L89").Value = 12162.5
$ws.Range("N89").Value = -23394.5

$ws.Range("H99").Value = 1252164
$ws.Range("I99").Value = 3334132.8
$ws.Range("J99").Value = 2982.8
$ws.Range("K99").Value = 3334132.8
$ws.Range("L99").Value = 2982.8
$ws.Range("M99").Value = -3332634.8
$ws.Range("N99").Value = -5978.8

$ws.Range("H126").Value = 1252164
$ws.Range("I126").Value = 3334132.8
$ws.Range("J126").Value = 2982.8
$ws.Range("K126").Value = 10002398.4
$ws.Range("L126").Value = 8948.400000000001
$ws.Range("M126").Value = -9999928.399999999
$ws.Range("N126").Value = -13888.4

$ws.Range("H136").Value = 2072000
$ws.Range("I136").Value = 2899858.5
$ws.Range("J136").Value = 2353.8333
$ws.Range("K136").Value = 8699575.5
$ws.Range("L136").Value = 7061.499899999999
$ws.Range("M136").Value = -8697025.5
$ws.Range("N136").Value = -12161.4999

$ws.Range("H140").Value = 58998.5
$ws.Range("J140").Value = 58998.5
$ws.Range("L140").Value = 58998.5
$ws.Range("N140").Value = -69358.5


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 924.75
$ws.Range("I32").Value = 849.5
$ws.Range("J32").Value = 1000
$ws.Range("K32").Value = 2548.5
$ws.Range("L32").Value = 3000
$ws.Range("M32").Value = -2265.5
$ws.Range("N32").Value = -3566

$ws.Range("H105").Value = 2760.6155
$ws.Range("J105").Value = 2916.6667
$ws.Range("L105").Value = 8750.000100000001
$ws.Range("N105").Value = -13992.0001

$ws.Range("H129").Value = 24181.709
$ws.Range("I129").Value = 388.4
$ws.Range("K129").Value = 1165.2
$ws.Range("M129").Value = 3834.8

$ws.Range("H131").Value = 837.72
$ws.Range("J131").Value = 844.04126
$ws.Range("L131").Value = 2532.12378
$ws.Range("N131").Value = -12612.12378

$ws.Range("H132").Value = 700
$ws.Range("I132").Value = 700
$ws.Range("K132").Value = 6300
$ws.Range("M132").Value = -3770


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 5365589
$ws.Range("I12").Value = 6364091
$ws.Range("K12").Value = 6364091
$ws.Range("M12").Value = -6363951

$ws.Range("H70").Value = 4475.643
$ws.Range("I70").Value = 4579
$ws.Range("J70").Value = 4398.125
$ws.Range("K70").Value = 4579
$ws.Range("L70").Value = 4398.125
$ws.Range("M70").Value = -4309
$ws.Range("N70").Value = -4938.125

$ws.Range("H73").Value = 4475.643
$ws.Range("I73").Value = 4579
$ws.Range("J73").Value = 4398.125
$ws.Range("K73").Value = 4579
$ws.Range("L73").Value = 4398.125
$ws.Range("M73").Value = -3643
$ws.Range("N73").Value = -6270.125

$ws.Range("H132").Value = 3499128.5
$ws.Range("I132").Value = 6411939.5
$ws.Range("K132").Value = 19235818.5
$ws.Range("M132").Value = -19233288.5

$ws.Range("H138").Value = 47059.668
$ws.Range("J138").Value = 47059.668
$ws.Range("L138").Value = 47059.668
$ws.Range("N138").Value = -57339.668


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3376.2222
$ws.Range("I7").Value = 2548.25
$ws.Range("K7").Value = 2548.25
$ws.Range("M7").Value = -2436.25

$ws.Range("H22").Value = 2201.75
$ws.Range("I22").Value = 2565
$ws.Range("J22").Value = 1942.2858
$ws.Range("K22").Value = 2565
$ws.Range("L22").Value = 1942.2858
$ws.Range("M22").Value = -2270
$ws.Range("N22").Value = -2532.2858

$ws.Range("H27").Value = 2201.75
$ws.Range("I27").Value = 2565
$ws.Range("J27").Value = 1942.2858
$ws.Range("K27").Value = 2565
$ws.Range("L27").Value = 1942.2858
$ws.Range("M27").Value = -2458
$ws.Range("N27").Value = -2156.2858

$ws.Range("H32").Value = 6183.5
$ws.Range("I32").Value = 4954.6665
$ws.Range("K32").Value = 4954.6665
$ws.Range("M32").Value = -4637.6665

$ws.Range("H40").Value = 3540.2307
$ws.Range("I40").Value = 1456.6364
$ws.Range("K40").Value = 1456.6364
$ws.Range("M40").Value = -1320.6364

$ws.Range("H46").Value = 2332.182
$ws.Range("I46").Value = 1500
$ws.Range("K46").Value = 1500
$ws.Range("M46").Value = -1312

$ws.Range("H122").Value = 10999.8
$ws.Range("J122").Value = 11666.667
$ws.Range("L122").Value = 35000.001
$ws.Range("N122").Value = -39900.001

$ws.Range("H126").Value = 3376.2222
$ws.Range("I126").Value = 2548.25
$ws.Range("K126").Value = 7644.75
$ws.Range("M126").Value = -5174.75

$ws.Range("H132").Value = 1857.0294
$ws.Range("I132").Value = 1760.5834
$ws.Range("J132").Value = 1909.6364
$ws.Range("K132").Value = 5281.7502
$ws.Range("L132").Value = 5728.9092
$ws.Range("M132").Value = -2751.7502
$ws.Range("N132").Value = -10788.9092

$ws.Range("H136").Value = 3265.5789
$ws.Range("I136").Value = 1962.25
$ws.Range("J136").Value = 5499.857
$ws.Range("K136").Value = 5886.75
$ws.Range("L136").Value = 16499.571
$ws.Range("M136").Value = -3336.75
$ws.Range("N136").Value = -21599.571


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 27000
$ws.Range("J92").Value = 29000
$ws.Range("L92").Value = 29000
$ws.Range("N92").Value = -33992

$ws.Range("H132").Value = 1307.6957
$ws.Range("I132").Value = 1027.5238
$ws.Range("K132").Value = 3082.5714
$ws.Range("M132").Value = -552.5713999999998

$ws.Range("H135").Value = 84166.5
$ws.Range("J135").Value = 84166.5
$ws.Range("L135").Value = 84166.5
$ws.Range("N135").Value = -94306.5

$ws.Range("H136").Value = 2947.8096
$ws.Range("I136").Value = 3918.6667
$ws.Range("J136").Value = 2219.6667
$ws.Range("K136").Value = 11756.0001
$ws.Range("L136").Value = 6659.000100000001
$ws.Range("M136").Value = -9206.000100000001
$ws.Range("N136").Value = -11759.0001

